$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold HYPERLINK(...) formulas needing a visible link-text
# second argument added: S, T, V, W, X, Y (U - "Knärotsbufferlänk" - has no
# formulas in this sheet).
$linkCols = @(19, 20, 22, 23, 24, 25)

$lastRow = 155

for ($r = 2; $r -le $lastRow; $r++) {

    # Column C ("Förändrad") moves from 45184 to 45186 for every data row.
    $cVal = $ws.Cells.Item($r, 3).Value2
    if ($cVal -eq 45184) {
        $ws.Cells.Item($r, 3).Value = 45186
    }

    # Rows with species hits also carry HYPERLINK formulas in S/T/V/W/X/Y;
    # append the case identifier (column A) as the link's friendly text.
    $ident = $ws.Cells.Item($r, 1).Value2

    foreach ($col in $linkCols) {
        $cell = $ws.Cells.Item($r, $col)
        $f = $cell.Formula
        if ($f -ne "" -and $f -like '*HYPERLINK(*' -and $f -notlike '*,*') {
            $newFormula = $f -replace '\)$', (', "' + $ident + '")')
            $cell.Formula = $newFormula
        }
    }
}
